$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab and update the cumulative-cutoff label text
$ws.Name = "Through 2022-06-19"
$ws.Range("B1").Value = "June 2022 (through June 19)"

# Updated / new counts for the "through June 19" dataset.
# Englewood
$ws.Range("B2").Value = 6
# Auburn Gresham
$ws.Range("AF3").Value = 1
# North Lawndale
$ws.Range("B4").Value = 6
# Humboldt Park
$ws.Range("H6").Value = 4
# Grand Crossing
$ws.Range("B9").Value = 5
$ws.Range("H9").Value = 3
# Roseland
$ws.Range("N12").Value = 3
# Austin
$ws.Range("B14").Value = 7
$ws.Range("AF14").Value = 3
# West Loop
$ws.Range("B17").Value = 2
# Little Italy, UIC
$ws.Range("H19").Value = 1
# Little Village
$ws.Range("B23").Value = 2
$ws.Range("H23").Value = 2
# Chicago Lawn
$ws.Range("AL24").Value = 1
# Morgan Park
$ws.Range("T35").Value = 1
# Chatham
$ws.Range("N38").Value = 3
# Calumet Heights
$ws.Range("T51").Value = 1
# Gage Park
$ws.Range("T57").Value = 1
# Galewood
$ws.Range("N58").Value = 1
# Old Town
$ws.Range("H81").Value = 1
# River North
$ws.Range("H85").Value = 2
$ws.Range("AF85").Value = 1
# Wicker Park
$ws.Range("N96").Value = 2
